# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# This script re-orders a handful of match rows within the "France National"
# sheet. For several small clusters of adjacent rows, the row identifier in
# column A stays put (it's just a running sequence number) while the rest of
# the row's data (match id, teams, odds, results, etc. in columns B..AD)
# moves to a different row in the same cluster. This mirrors a re-sort of
# the underlying source data that happened upstream, while keeping the
# sequential numbering in column A intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-match data (A = running index and is never touched;
# C = constant "France National"; D = match date, identical for every row
# inside each cluster below, so it is also left alone).
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Maps: destination row -> source row (i.e. destination row ends up holding
# the data that currently lives in the source row).
$rowMap = [ordered]@{
    18 = 19
    19 = 18
    21 = 23
    22 = 21
    23 = 22
    257 = 258
    258 = 259
    259 = 257
    260 = 261
    261 = 262
    262 = 260
    304 = 305
    305 = 306
    306 = 307
    307 = 308
    308 = 309
    309 = 304
}

# Union of every row referenced above (as both source and destination) --
# snapshot all of their current values BEFORE writing anything, since some
# rows are both a source for one destination and a destination for another.
$rowsToSnapshot = New-Object System.Collections.Generic.HashSet[int]
foreach ($key in $rowMap.Keys) {
    [void]$rowsToSnapshot.Add([int]$key)
    [void]$rowsToSnapshot.Add([int]$rowMap[$key])
}

$snapshot = @{}
foreach ($r in $rowsToSnapshot) {
    $rowValues = @{}
    foreach ($c in $cols) {
        $rowValues[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowValues
}

# Now apply the new values using the snapshot taken above (never re-reading
# from the sheet so earlier writes in this loop cannot clobber later reads).
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcValues[$c]
    }
}
